# Trade #21 closed at 2026-02-16 21:25:35 - momentum DOWN +0.000%
#
# Story: leadlag trade #3 (row 4 on the "leadlag" sheet) gets closed out
# with a +0.4354% / +$4.35 win via a 5-minute time exit, which is also
# appended as a row to "All Trades". A brand-new momentum trade #21 is
# opened and appended to the "momentum" sheet. The "Summary" and
# "Comparison" roll-up sheets are refreshed to reflect the new trade
# counts / win rates / P&L stats.
#
# NOTE: several of the roll-up values are stored as *literal text*
# (e.g. "33.3%", "0.92") rather than numbers, matching the source
# workbook's convention. Plain `Range.Value = "33.3%"` gets silently
# reinterpreted by Excel as a percentage number (and "0.92"/"1.84" as
# plain numbers), so those assignments use a leading apostrophe to force
# literal text, then restore the "Normal" style so no stray per-cell
# number-format/quote-prefix is left behind.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet: update OVERALL + leadlag STRATEGY rows
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("C2").Value = 3

$wsSummary.Range("D2").Value = "'33.3%"
$wsSummary.Range("D2").Style = "Normal"
$wsSummary.Range("E2").Value = "'-0.0385%"
$wsSummary.Range("E2").Style = "Normal"
$wsSummary.Range("F2").Value = "'-0.0128%"
$wsSummary.Range("F2").Style = "Normal"

$wsSummary.Range("C3").Value = 18

$wsSummary.Range("D3").Value = "'5.6%"
$wsSummary.Range("D3").Style = "Normal"
$wsSummary.Range("E3").Value = "'-0.0385%"
$wsSummary.Range("E3").Style = "Normal"
$wsSummary.Range("F3").Value = "'-0.0021%"
$wsSummary.Range("F3").Style = "Normal"

# ---------------------------------------------------------------------
# 2. leadlag sheet: close out trade #3 (row 4)
# ---------------------------------------------------------------------
$wsLeadlag = $wb.Worksheets.Item("leadlag")
$wsLeadlag.Range("G4").Value = 69173.562481
$wsLeadlag.Range("H4").Value = "CLOSED"
$wsLeadlag.Range("I4").Value = 0.4354
$wsLeadlag.Range("J4").Value = 4.35
$wsLeadlag.Range("M4").Value = "time_exit_5min"
$wsLeadlag.Range("N4").Value = 5

# ---------------------------------------------------------------------
# 3. momentum sheet: append new trade #21 (row 4)
# ---------------------------------------------------------------------
$wsMomentum = $wb.Worksheets.Item("momentum")
$wsMomentum.Range("A4").Value = 21
$wsMomentum.Range("B4").Value = "'2026-02-16"
$wsMomentum.Range("B4").Style = "Normal"
$wsMomentum.Range("C4").Value = "21:25:35"
$wsMomentum.Range("D4").Value = "momentum"
$wsMomentum.Range("E4").Value = "DOWN"
$wsMomentum.Range("F4").Value = 69077.44500000001
$wsMomentum.Range("H4").Value = "OPEN"
$wsMomentum.Range("I4").Value = 0
$wsMomentum.Range("J4").Value = 0
$wsMomentum.Range("K4").Value = 0.9
$wsMomentum.Range("L4").Value = "Downward momentum: -0.283% over 10 samples"
$wsMomentum.Range("N4").Value = 0

# ---------------------------------------------------------------------
# 4. All Trades sheet: append closed trade #3 (row 4)
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")
$wsAll.Range("A4").Value = 3
$wsAll.Range("B4").Value = "'2026-02-16"
$wsAll.Range("B4").Style = "Normal"
$wsAll.Range("C4").Value = "21:20:31"
$wsAll.Range("D4").Value = "leadlag"
$wsAll.Range("E4").Value = "DOWN"
$wsAll.Range("F4").Value = 69476.05
$wsAll.Range("G4").Value = 69173.562481
$wsAll.Range("H4").Value = "CLOSED"
$wsAll.Range("I4").Value = 0.4354
$wsAll.Range("J4").Value = 4.35
$wsAll.Range("K4").Value = 0.75
$wsAll.Range("L4").Value = "Binance leading with -0.099% move"
$wsAll.Range("M4").Value = "time_exit_5min"
$wsAll.Range("N4").Value = 5

# ---------------------------------------------------------------------
# 5. Comparison sheet: update leadlag row
# ---------------------------------------------------------------------
$wsComparison = $wb.Worksheets.Item("Comparison")
$wsComparison.Range("B2").Value = 18

$wsComparison.Range("C2").Value = "'5.6%"
$wsComparison.Range("C2").Style = "Normal"
$wsComparison.Range("D2").Value = "'0.92"
$wsComparison.Range("D2").Style = "Normal"
$wsComparison.Range("E2").Value = "'+0.4354%"
$wsComparison.Range("E2").Style = "Normal"
$wsComparison.Range("G2").Value = "'1.84"
$wsComparison.Range("G2").Style = "Normal"
